$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.078.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.89%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.273.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'185.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.42%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'576.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.59%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.18%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.73%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.18%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.408"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.11%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'3.841.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.20%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -0.50%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'27.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.10%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'67.369.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.53%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -1.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.268.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.49%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'440.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +9.86%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.02%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'13.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.57%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +2.18%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'74.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.38%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.26%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'Polygon"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.512"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.15%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "'WrappedeETH"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").Value = "'3.433.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.64%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.46%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  +0.18%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -4.56%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.06%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +1.04%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'22.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.29%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'5.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.65%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -0.01%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -0.77%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.05%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'1.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +5.09%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'162.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.15%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'Stacks"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'1.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.41%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'EnergySwap"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'27.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.59%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.62%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -2.95%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'2.708.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.18%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -0.84%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'40.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.02%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -1.16%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'24.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.57%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -2.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'326.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.34%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.02%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'31.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +2.83%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +2.07%  "
$ws.Range("E51").Style = "Normal"
